$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header row: add Wins, Losses, Ties columns after the existing last column (AC)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style (bold / bordered / centered) from an existing header cell
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the season record values for every data row (2 through 60)
for ($r = 2; $r -le 60; $r++) {
    $ws.Cells.Item($r, 30).Value = 99
    $ws.Cells.Item($r, 31).Value = 63
    $ws.Cells.Item($r, 32).Value = 0
}
